# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text even when it looks like a number
# (keeps "1.00"/"0.370"-style trailing zeros intact instead of letting Excel
# coerce the input to a numeric Value), then restores the default cell style
# so we do not leave a stray text-format override behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '62.609.19'
$ws.Range('E2').Value = '  -1.94%  '
$ws.Range('D3').Value = '2.553.07'
$ws.Range('E3').Value = '  -2.31%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue $ws.Range('D5') '562.55'
$ws.Range('E5').Value = '  -2.29%  '
Set-TextValue $ws.Range('D6') '151.63'
$ws.Range('E6').Value = '  -3.29%  '
Set-TextValue $ws.Range('D7') '1.00'
$ws.Range('E7').Value = '  +0.08%  '
Set-TextValue $ws.Range('D8') '0.612'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = '2.557.98'
$ws.Range('E9').Value = '  -2.01%  '
Set-TextValue $ws.Range('D10') '0.112'
$ws.Range('E10').Value = '  -5.94%  '
Set-TextValue $ws.Range('D11') '5.58'
$ws.Range('E11').Value = '  -4.32%  '
$ws.Range('E12').Value = '  -0.29%  '
Set-TextValue $ws.Range('D13') '0.370'
$ws.Range('E13').Value = '  -3.43%  '
Set-TextValue $ws.Range('D14') '27.26'
$ws.Range('E14').Value = '  -3.55%  '
$ws.Range('D15').Value = '3.022.75'
$ws.Range('E15').Value = '  -2.11%  '
Set-TextValue $ws.Range('D16') '0.0000174'
$ws.Range('E16').Value = '  -4.89%  '
$ws.Range('D17').Value = '62.506.72'
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('D18').Value = '2.578.81'
$ws.Range('E18').Value = '  -1.15%  '
Set-TextValue $ws.Range('D19') '11.66'
$ws.Range('E19').Value = '  -3.23%  '
Set-TextValue $ws.Range('D20') '7.28'
$ws.Range('E20').Value = '  -4.71%  '
Set-TextValue $ws.Range('D21') '4.36'
$ws.Range('E21').Value = '  -4.71%  '
Set-TextValue $ws.Range('D22') '331.12'
$ws.Range('E22').Value = '  -3.49%  '
$ws.Range('E23').Value = '  -0.07%  '
Set-TextValue $ws.Range('D24') '66.42'
$ws.Range('E24').Value = '  -1.53%  '
Set-TextValue $ws.Range('D25') '1.81'
$ws.Range('E25').Value = '  +3.27%  '
Set-TextValue $ws.Range('D26') '0.0000107'
$ws.Range('E26').Value = '  -2.95%  '
$ws.Range('D27').Value = '2.698.44'
$ws.Range('E27').Value = '  -1.83%  '
Set-TextValue $ws.Range('D28') '8.86'
$ws.Range('E28').Value = '  -3.65%  '
Set-TextValue $ws.Range('D29') '1.55'
$ws.Range('E29').Value = '  -0.98%  '
Set-TextValue $ws.Range('D30') '548.31'
$ws.Range('E30').Value = '  -7.73%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range('D31') '1.00'
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D32') '7.85'
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('E33').Value = '  -2.39%  '
Set-TextValue $ws.Range('D34') '1.97'
$ws.Range('E34').Value = '  -4.87%  '
Set-TextValue $ws.Range('D35') '1.65'
$ws.Range('E35').Value = '  -6.50%  '
Set-TextValue $ws.Range('D36') '6.26'
$ws.Range('E36').Value = '  -4.99%  '
Set-TextValue $ws.Range('D37') '5.11'
$ws.Range('E37').Value = '  -4.90%  '
Set-TextValue $ws.Range('D38') '1.00'
$ws.Range('E38').Value = '  +0.22%  '
Set-TextValue $ws.Range('D39') '0.389'
$ws.Range('E39').Value = '  -4.51%  '
Set-TextValue $ws.Range('D40') '19.15'
$ws.Range('E40').Value = '  -3.07%  '
Set-TextValue $ws.Range('D41') '152.58'
$ws.Range('E41').Value = '  -1.01%  '
Set-TextValue $ws.Range('D42') '1.80'
$ws.Range('E42').Value = '  -3.91%  '
Set-TextValue $ws.Range('D43') '1.00'
$ws.Range('E43').Value = '  +0.06%  '
Set-TextValue $ws.Range('D44') '2.38'
$ws.Range('E44').Value = '  -2.93%  '
Set-TextValue $ws.Range('D45') '153.92'
$ws.Range('E45').Value = '  -1.05%  '
Set-TextValue $ws.Range('D46') '22.82'
$ws.Range('E46').Value = '  -0.75%  '
Set-TextValue $ws.Range('D47') '3.75'
$ws.Range('E47').Value = '  -4.24%  '
Set-TextValue $ws.Range('D48') '0.0563'
$ws.Range('E48').Value = '  -4.91%  '
Set-TextValue $ws.Range('D49') '0.615'
$ws.Range('E49').Value = '  -2.11%  '
Set-TextValue $ws.Range('D50') '0.0975'
$ws.Range('E50').Value = '  -4.08%  '
Set-TextValue $ws.Range('D51') '0.0240'
$ws.Range('E51').Value = '  -3.00%  '
